# Thesis proposal edit: extend the "streaming data" discussion with three
# more sentences, remove the _GoBack bookmark from its old spot, then add
# the whole "IDR/QR" section (heading + three body paragraphs) below a
# run of blank paragraphs (page-break filler), finally re-homing the
# _GoBack bookmark at the end of the "common way" paragraph.

$d = $word.ActiveDocument

# --- 1. Locate + remove the existing (hidden) _GoBack bookmark -----------
# It currently sits right after "...different algorithms " at the end of
# the last paragraph; it needs to move further down in the new text.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# --- 2. Append the three new runs to the end of that same paragraph ------
$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertAfter("which work fine for static data but not perform well for the streaming data. For example, principal component analysis (PCA), Linear Disc")
$tail.Collapse(0)
$tail.InsertAfter("riminant analysis (LDA), Maximum Marginal classifier (MMC). Previous study show that these are not well ")
$tail.Collapse(0)
$tail.InsertAfter("suitable for the ")

# --- 3. Append the rest of the new content as a block of OOXML -----------
$endRange = $d.Content
$endRange.Collapse(0)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>IDR/QR: An Incremental Dimension Reduction Algorithm via QR Decomposition</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Dimension reduction are necessary for many database and data mining application mainly for efficient storage and retrieval of high dimensional data.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> For larger data set it will be better to not store whole data matrix in the memory. More importantly</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> when the new data items arrive the algorithm should c</w:t>
  </w:r>
  <w:r>
    <w:t>onstrain the computational cost</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> The almost all effort for having efficient storage deteriorates with the increment of the dimensions. The solution of this dimensional curse is to reduce the dimension by some means and then apply the multi indexing techniques.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">The common way </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">to reduce computationally and time is to have a small chunk of data available over a certain fixed </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>period of time</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t>.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>One of the main problem find in almost all paper is: it is difficult to design an incremental solution for the eigenvalue problem on the product of scatter matrixes.</w:t>
  </w:r>
</w:p>
'@

$endRange.InsertXML($newXml)
